$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlUp = -4162
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End($xlUp).Row

# ---------------------------------------------------------------------------
# 1) Row 3 / Row 4 content update.
#    A new inspection ("A 13111-2022") ends up in row 3 with refreshed stats
#    (one more signal species found, one more red-listed species, "Svavelriska"
#    added to the species list). The inspection that used to sit in row 3
#    ("A 11731-2019") moves down to row 4 unchanged (its "Förändrad" date is
#    refreshed along with every other row further below).
# ---------------------------------------------------------------------------

# New row 3: "A 13111-2022"
$ws.Cells.Item(3, 1).Value = "A 13111-2022"             # Beteckning
$ws.Cells.Item(3, 2).Value = 44644                       # Datum
$ws.Cells.Item(3, 3).Value = 45186                       # Förändrad
$ws.Cells.Item(3, 4).Value = "VÄSTRA GÖTALANDS LÄN"      # Län
$ws.Cells.Item(3, 5).Value = "ALE"                       # Kommun
$ws.Cells.Item(3, 7).Value = 10.3                        # Area (ha)
$ws.Cells.Item(3, 8).Value = 4                           # Fridlysta
$ws.Cells.Item(3, 9).Value = 11                          # Signalarter
$ws.Cells.Item(3, 10).Value = 4                          # NT
$ws.Cells.Item(3, 11).Value = 0                          # VU
$ws.Cells.Item(3, 12).Value = 0                          # EN
$ws.Cells.Item(3, 13).Value = 0                          # CR
$ws.Cells.Item(3, 14).Value = 0                          # RE
$ws.Cells.Item(3, 15).Value = 4                           # Rödlistade
$ws.Cells.Item(3, 16).Value = 0                           # Hotade
$ws.Cells.Item(3, 17).Value = 15                          # Alla arter
$ws.Cells.Item(3, 18).Value = "Entita`r`nMindre hackspett`r`nSpillkråka`r`nTalltita`r`nBlåmossa`r`nBronshjon`r`nFällmossa`r`nGuldlockmossa`r`nKlippfrullania`r`nRödgul trumpetsvamp`r`nSmal svampklubba`r`nStor revmossa`r`nSvavelriska`r`nVästlig hakmossa`r`nVågbandad barkbock"

$ws.Cells.Item(3, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/artfynd/A 13111-2022.xlsx", "A 13111-2022")'
$ws.Cells.Item(3, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/kartor/A 13111-2022.png", "A 13111-2022")'
$ws.Cells.Item(3, 21).Value = ""
$ws.Cells.Item(3, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/klagomål/A 13111-2022.docx", "A 13111-2022")'
$ws.Cells.Item(3, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/klagomålsmail/A 13111-2022.docx", "A 13111-2022")'
$ws.Cells.Item(3, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/tillsyn/A 13111-2022.docx", "A 13111-2022")'
$ws.Cells.Item(3, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/tillsynsmail/A 13111-2022.docx", "A 13111-2022")'

# New row 4: "A 11731-2019" (same figures it always had, only the
# "Förändrad" date moves forward with the rest of the sheet)
$ws.Cells.Item(4, 1).Value = "A 11731-2019"
$ws.Cells.Item(4, 2).Value = 43518
$ws.Cells.Item(4, 3).Value = 45186
$ws.Cells.Item(4, 4).Value = "VÄSTRA GÖTALANDS LÄN"
$ws.Cells.Item(4, 5).Value = "ALE"
$ws.Cells.Item(4, 7).Value = 1.9
$ws.Cells.Item(4, 8).Value = 3
$ws.Cells.Item(4, 9).Value = 11
$ws.Cells.Item(4, 10).Value = 2
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 2
$ws.Cells.Item(4, 16).Value = 0
$ws.Cells.Item(4, 17).Value = 14
$ws.Cells.Item(4, 18).Value = "Kråka`r`nMindre hackspett`r`nBlåsfliksmossa`r`nBrandticka`r`nGrov fjädermossa`r`nHavstulpanlav`r`nKlippfrullania`r`nKornknutmossa`r`nSmal svampklubba`r`nSotriska`r`nStor revmossa`r`nStubbspretmossa`r`nVästlig hakmossa`r`nVanlig padda"

$ws.Cells.Item(4, 19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/artfynd/A 11731-2019.xlsx", "A 11731-2019")'
$ws.Cells.Item(4, 20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/kartor/A 11731-2019.png", "A 11731-2019")'
$ws.Cells.Item(4, 21).Value = ""
$ws.Cells.Item(4, 22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/klagomål/A 11731-2019.docx", "A 11731-2019")'
$ws.Cells.Item(4, 23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/klagomålsmail/A 11731-2019.docx", "A 11731-2019")'
$ws.Cells.Item(4, 24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/tillsyn/A 11731-2019.docx", "A 11731-2019")'
$ws.Cells.Item(4, 25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_ALE/tillsynsmail/A 11731-2019.docx", "A 11731-2019")'

# Re-writing the wrapped, multi-line species cell (column R) makes Excel
# auto-fit the row to the new text; put the fixed row height back the way
# every other row in the sheet has it.
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(4).RowHeight = 15

# ---------------------------------------------------------------------------
# 2) Every row's "Förändrad" date (column C) moves from 2023-09-20 (45184)
#    to 2023-09-22 (45186). Rows 3 & 4 are already handled above.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 3 -or $r -eq 4) { continue }
    $c = $ws.Cells.Item($r, 3)
    if ($c.Value2 -eq 45184) {
        $c.Value = 45186
    }
}

# ---------------------------------------------------------------------------
# 3) Every existing HYPERLINK(...) formula (columns S, T, U, V, W, X, Y) gets
#    a second argument added: the friendly display text, which is simply the
#    "Beteckning" (case id) of that row.
# ---------------------------------------------------------------------------
$hyperlinkCols = 19, 20, 21, 22, 23, 24, 25  # S,T,U,V,W,X,Y
for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 3 -or $r -eq 4) { continue }
    $beteckning = $ws.Cells.Item($r, 1).Value2
    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Cells.Item($r, $col)
        $f = $cell.Formula
        if ($f -and $f.StartsWith("=HYPERLINK(") -and $f.EndsWith('")') -and -not $f.Contains(', "')) {
            $newf = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newf
        }
    }
}
